# Adds results (read counts) for the "Combined" block (rows 31-35) and the
# "Mock-2" block (rows 49-53), matching the pattern already present for the
# earlier blocks in the same worksheet (e.g. rows 13-17, 37-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Combined block (rows 31-35) ---
$ws.Range("E31").Value = 10881
$ws.Range("G31").Value = 10860

$ws.Range("E32").Value = 10840
$ws.Range("G32").Value = 10794

$ws.Range("F33").Formula = "=SUM(E31:E32)"
$ws.Range("G33").Formula = "=SUM(G31:G32)"

$ws.Range("E34").Value = 0
$ws.Range("E35").Value = 0

# --- Mock-2 block (rows 49-53) ---
$ws.Range("E49").Value = 10877
$ws.Range("G49").Value = 10861

$ws.Range("E50").Value = 10866
$ws.Range("G50").Value = 10830

$ws.Range("F51").Formula = "=SUM(E49:E50)"
$ws.Range("G51").Formula = "=SUM(G49:G50)"

$ws.Range("E52").Value = 0
$ws.Range("E53").Value = 0

# --- Update the saved selection to match the final cursor position ---
$ws.Range("G52").Select()
